# Auto-generated edit script: updates crypto price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as plain text so numeric-looking
# strings (e.g. '1.00', '14.30', '0.0225') are not coerced into numbers
# and lose their original formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.071.61'
$ws.Range("E2").Value = '  +1.52%  '
$ws.Range("D3").Value = '2.062.30'
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '249.42'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = '0.672'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '55.43'
$ws.Range("E8").Value = '  +16.45%  '
$ws.Range("D9").Value = '60.99'
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("E11").Value = '  +6.12%  '
$ws.Range("E12").Value = '  +6.03%  '
$ws.Range("D13").Value = '15.04'
$ws.Range("E13").Value = '  +3.41%  '
$ws.Range("D14").Value = '2.359.71'
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("D15").Value = '0.819'
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").Value = '5.31'
$ws.Range("E16").Value = '  +4.47%  '
$ws.Range("D17").Value = '2.062.23'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").Value = '37.025.62'
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").Value = '0.0₃0928'
$ws.Range("E19").Value = '  +11.83%  '
$ws.Range("D20").Value = '73.43'
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("D21").Value = '14.30'
$ws.Range("E21").Value = '  +8.25%  '
$ws.Range("D22").Value = '5.37'
$ws.Range("E22").Value = '  +2.45%  '
$ws.Range("D23").Value = '237.66'
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  -2.66%  '
$ws.Range("D26").Value = '170.69'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("D28").Value = '20.20'
$ws.Range("E28").Value = '  -5.44%  '
$ws.Range("D29").Value = '1.99'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("E30").Value = '  +1.59%  '
$ws.Range("D31").Value = '4.61'
$ws.Range("E31").Value = '  +2.79%  '
$ws.Range("E32").Value = '  +6.63%  '
$ws.Range("D33").Value = '0.0627'
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("E34").Value = '  +7.10%  '
$ws.Range("D35").Value = '0.0886'
$ws.Range("E35").Value = '  -2.91%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '2.29'
$ws.Range("E37").Value = '  -5.72%  '
$ws.Range("E38").Value = '  -4.47%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("E40").Value = '  +22.23%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0225'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '17.72'
$ws.Range("E42").Value = '  +11.16%  '
$ws.Range("E43").Value = '  -1.82%  '
$ws.Range("D44").Value = '97.02'
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("B46").Value = 'Gas'
$ws.Range("C46").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D46").Value = '14.10'
$ws.Range("E46").Value = '  -51.06%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.07'
$ws.Range("E47").Value = '  +47.13%  '
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").Value = '  +6.99%  '
$ws.Range("D49").Value = '1.298.80'
$ws.Range("E49").Value = '  -2.29%  '
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").Value = '4.15'
$ws.Range("E51").Value = '  +8.83%  '
